# Apply "user stories spreadsheet updated to reflect completion of payment
# system" edits to the User_Stories sheet of the workbook.
#
# Summary of the change:
#  - The "Purchasing & Checkout" items that cover the shopping-bag / checkout
#    flow (rows 22-26) are marked complete with a new status string
#    "Complete  - Bag app".
#  - User story #22 (Review & Rating, row 28) and a brand-new user story #26
#    (row 33, "Edit Shipping cost per country/continent") are highlighted
#    with a new light-blue fill.
#  - The blank spacer row that used to sit above the "Admin & Product
#    Management" section header is removed - the header itself moves up one
#    row (now row 29) and a new Admin row (#23 "Add products") takes its old
#    slot (row 30). Every following Admin row shifts its User Story ID down
#    by one to make room for the newly appended story #26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User_Stories")

# ---------------------------------------------------------------------
# 1) Rows 22-26: "Complete  - Bag app" status (note: two spaces before the
#    hyphen, matching the author's original text exactly).
# ---------------------------------------------------------------------
$bagAppStatus = "Complete  - Bag app"
$ws.Range("E22").Value = $bagAppStatus
$ws.Range("E23").Value = $bagAppStatus
$ws.Range("E24").Value = $bagAppStatus
$ws.Range("E25").Value = $bagAppStatus
$ws.Range("E26").Value = $bagAppStatus

# ---------------------------------------------------------------------
# 2) Highlight row 28 (User story 22 - Review & Rating) with the new
#    light-blue fill (RGB 0,176,240 => 0x00F0B000 in COM's BGR order).
# ---------------------------------------------------------------------
$highlightColor = 15773696  # RGB(0, 176, 240)
$ws.Range("A28:E28").Interior.Color = $highlightColor

# ---------------------------------------------------------------------
# 3) Remove the blank spacer row (old row 29) by turning it into the
#    "Admin & Product Management" section header (formerly row 30), then
#    reuse the vacated old-header row (row 30) as a normal data row.
# ---------------------------------------------------------------------

# 3a. Un-merge the old header band before reshuffling its contents.
$ws.Range("A30:E30").UnMerge()

# 3b. Row 29 becomes the section header: copy the header formatting from an
#     existing section header (row 21, "Purchasing & Checkout") so the fonts
#     / fills / borders / center alignment match exactly, then merge + set
#     the header text.
$ws.Range("A21:E21").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A29:E29").Merge()
$ws.Range("A29").Value = "Admin & Product Management"

# 3c. Row 30 becomes a plain data row (User story 23 - "Add products"): copy
#     the plain-row formatting from row 27 (a normal, unhighlighted data
#     row), then fill in the values.
$ws.Range("A27:E27").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A30").Value = 23
$ws.Range("B30").Value = "Admin"
$ws.Range("C30").Value = "Add products"
$ws.Range("D30").Value = "Add new products to site"
$ws.Range("E30").ClearContents()

# ---------------------------------------------------------------------
# 4) Shift the remaining Admin rows' User Story IDs down by one (23->24,
#    24->25, 25->26) to make room for the newly appended story.
# ---------------------------------------------------------------------
$ws.Range("A31").Value = 24
$ws.Range("B31").Value = "Admin"
$ws.Range("C31").Value = "Edit/update existing products"
$ws.Range("D31").Value = "Change details of existing products"

$ws.Range("A32").Value = 25
$ws.Range("B32").Value = "Admin"
$ws.Range("C32").Value = "Delete products"
$ws.Range("D32").Value = "Remove products no longer for sale on the site"

# ---------------------------------------------------------------------
# 5) New row 33: User story 26 - "Edit Shipping cost per country/continent",
#    highlighted the same way as row 28.
# ---------------------------------------------------------------------
$ws.Range("A33").Value = 26
$ws.Range("B33").Value = "Admin"
$ws.Range("C33").Value = "Edit Shipping cost per country/continent"
$ws.Range("D33").Value = "Option to change shipping cost depending on the location of the customer"
$ws.Range("E33").ClearContents()
$ws.Range("A33:E33").Interior.Color = $highlightColor

# ---------------------------------------------------------------------
# 6) Restore the active selection to match the author's final cursor
#    position.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E26").Select()
